{"js": "// Helper: find the first occurrence of `needle` in the document body and\n// replace it in place with `replacement` (preserving surrounding runs).\nasync function replaceText(body, needle, replacement) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replacement, \"Replace\");\n    await context.sync();\n  }\n  return results.items.length > 0;\n}\n\nconst body = context.document.body;\n\n// 1. Update the date.\nawait replaceText(\n  body,\n  \"2023-01-28\",\n  \"2023-01-31\"\n);\n\n// 2. Rewrite the \"focus\" paragraph.\nawait replaceText(\n  body,\n  \"With this role my focus has been on growing teams while retaining culture and aligning teams with purpose, values and business goals.\",\n  \"My focus has been to grow the Cloud team whilst creating and evolving a culture that is aligned to the purpose, values and goals for the business.\"\n);\n\n// 3. Rewrite the \"responsibilities\" paragraph.\nawait replaceText(\n  body,\n  \"My responsibilities continue to be in shaping the product portfolio as part of the business unit\\u2019s leadership team, ensure that my teams have the skills, the tools and the autonomy to delivery successfully to customers. This includes building a platform of tools to support the delivery, operations and team learning to ensure the practice can continue to expand from 50 to 100+ engineers.s\",\n  \"As part of the leadership team, I am responsible for shaping the product portfolio ensuring that my teams have the skills, the tools and the autonomy to deliver successfully to customers. This includes building a platform of tools to support the delivery, operations and team learning to ensure the practice can continue to expand from 50 to 100+ engineers.\"\n);\n\n// 4. Split the \"While continuing to manage\" paragraph into two paragraphs.\nawait replaceText(\n  body,\n  \"While continuing to manage by initial team of 8 engineers, with the technical principal role I was much more engaged across out multi-cloud portfolio (Azure, AWS and GCP) supporting the expansion and definition of new service offerings to customers. This involved continued end-customer engagement and consulting at multi-levels, working on business plans and service definitions. This role expanded my knowledge of Product Management, customer leadership and allowed to be start my journey of learning strategy and roadmap development.\",\n  \"This was a combined management and delivery role. I managed the initial team of 8 engineers, as the Technical Principal. I was much more engaged across the multi-cloud portfolio (Azure, AWS and GCP) supporting the expansion and definition of new service offerings to customers.\\nThis involved continued end-customer engagement and consulting at multi-levels, working on business plans and service definitions. This role expanded my knowledge of Product Management, customer leadership and allowed me to enhance my journey of strategy and roadmap development.\"\n);\n\n// 5. Windows capitalization + entrepreneurial/leadership wording.\n//    Two separate replacements so the run structure / the untouched\n//    \" \" run in between is preserved, matching the source edit.\nawait replaceText(\n  body,\n  \"With a windows team already established I took on a new role to build and grow the team and capability with Microsoft Azure.\",\n  \"With a Windows team already established I took on a new role to build and grow the team and capability with Microsoft Azure.\"\n);\nawait replaceText(\n  body,\n  \"This was a very entrepreneurial and leadership role within the business where I was leading all elements of the business growth including hiring the team, partner management with Microsoft and pre-sales with C-level customer prospects.\",\n  \"This was a very entrepreneurial/leadership role within the business where I was leading all elements of the business growth including hiring the team, partner management with Microsoft and pre-sales with C-level customer prospects.\"\n);\n\n// 6. \"I was also expanding\" -> \"I also expanded ... worked with my team\".\nawait replaceText(\n  body,\n  \"I was also expanding my technical knowledge of the Azure platform and working with my team to introduce new tools and establish a deployment and operations platform to support large customer engagements on this new (for the business) platform.\",\n  \"I also expanded my technical knowledge of the Azure platform and worked with my team to introduce new tools and establish a deployment and operations platform to support large customer engagements on this new (for the business) platform.\"\n);\n\n// 7. \"has involved\" -> \"involved\".\nawait replaceText(\n  body,\n  \"A significant portion of this role has involved me liaising with and training development and operations teams both in London and in San Francisco. It also required communicating with all layers of management to ensure the successful rollout of the projects.\",\n  \"A significant portion of this role involved me liaising with and training development and operations teams both in London and in San Francisco. It also required communicating with all layers of management to ensure the successful rollout of the projects.\"\n);\n\n// 8. \"code open source\" -> \"open-source code\"; \"I have found\" -> \"I found\".\nawait replaceText(\n  body,\n  \"This work gave me the opportunity to work and release a lot of code open source and build communities around that; in particular I made significant contributions to the Puppet community by providing additional support for Windows. I learnt a lot from this experience in building community though blogging, giving talks and podcasting on this and other subjects. Community evangelism is not something that is often thought of as significant in an operations-based role, but I have found it to be both enjoyable and beneficial on many occasions.\",\n  \"This work gave me the opportunity to work and release a lot of open-source code and build communities around that; in particular I made significant contributions to the Puppet community by providing additional support for Windows. I learnt a lot from this experience in building community though blogging, giving talks and podcasting on this and other subjects. Community evangelism is not something that is often thought of as significant in an operations-based role, but I found it to be both enjoyable and beneficial on many occasions.\"\n);\n\n// 9. \"I was leading\" -> \"I lead\"; \"we see value\" -> \"we saw value\".\nawait replaceText(\n  body,\n  \"From April 2015 to July 2016, I was leading my team in its projects and decision making. We remained autonomous within the organization, identifying problems and providing solutions wherever we see value.\",\n  \"From April 2015 to July 2016, I lead my team in its projects and decision making. We remained autonomous within the organization, identifying problems and providing solutions wherever we saw value.\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the date.\n$find = $d.Content.Find\n$find.Execute(\"2023-01-28\", $false, $false, $false, $false, $false, $true, 1, $false, \"2023-01-31\", 2)\n\n# 2. Rewrite the \"focus\" paragraph.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*With this role my focus has been on growing teams*\") {\n        $p.Range.Text = \"My focus has been to grow the Cloud team whilst creating and evolving a culture that is aligned to the purpose, values and goals for the business.\"\n        break\n    }\n}\n\n# 3. Rewrite the \"responsibilities\" paragraph.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*My responsibilities continue to be in shaping the product portfolio*\") {\n        $p.Range.Text = \"As part of the leadership team, I am responsible for shaping the product portfolio ensuring that my teams have the skills, the tools and the autonomy to deliver successfully to customers. This includes building a platform of tools to support the delivery, operations and team learning to ensure the practice can continue to expand from 50 to 100+ engineers.\"\n        break\n    }\n}\n\n# 4. Split the \"While continuing to manage\" paragraph into two paragraphs.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*While continuing to manage by initial team of 8 engineers*\") {\n        $p.Range.Text = \"This was a combined management and delivery role. I managed the initial team of 8 engineers, as the Technical Principal. I was much more engaged across the multi-cloud portfolio (Azure, AWS and GCP) supporting the expansion and definition of new service offerings to customers.`rThis involved continued end-customer engagement and consulting at multi-levels, working on business plans and service definitions. This role expanded my knowledge of Product Management, customer leadership and allowed me to enhance my journey of strategy and roadmap development.\"\n        break\n    }\n}\n\n# 5. Windows capitalization + entrepreneurial/leadership wording (same paragraph, two runs).\n$find = $d.Content.Find\n$find.Execute(\"With a windows team already established\", $false, $false, $false, $false, $false, $true, 1, $false, \"With a Windows team already established\", 2)\n\n$find = $d.Content.Find\n$find.Execute(\"This was a very entrepreneurial and leadership role\", $false, $false, $false, $false, $false, $true, 1, $false, \"This was a very entrepreneurial/leadership role\", 2)\n\n# 6. \"I was also expanding\" -> \"I also expanded ... worked with my team\".\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*I was also expanding my technical knowledge of the Azure platform*\") {\n        $p.Range.Text = \"I also expanded my technical knowledge of the Azure platform and worked with my team to introduce new tools and establish a deployment and operations platform to support large customer engagements on this new (for the business) platform.\"\n        break\n    }\n}\n\n# 7. \"has involved\" -> \"involved\".\n$find = $d.Content.Find\n$find.Execute(\"A significant portion of this role has involved me\", $false, $false, $false, $false, $false, $true, 1, $false, \"A significant portion of this role involved me\", 2)\n\n# 8. \"code open source\" -> \"open-source code\"; \"I have found\" -> \"I found\".\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*This work gave me the opportunity to work and release a lot of code open source*\") {\n        $p.Range.Text = \"This work gave me the opportunity to work and release a lot of open-source code and build communities around that; in particular I made significant contributions to the Puppet community by providing additional support for Windows. I learnt a lot from this experience in building community though blogging, giving talks and podcasting on this and other subjects. Community evangelism is not something that is often thought of as significant in an operations-based role, but I found it to be both enjoyable and beneficial on many occasions.\"\n        break\n    }\n}\n\n# 9. \"I was leading\" -> \"I lead\"; \"we see value\" -> \"we saw value\".\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*From April 2015 to July 2016, I was leading my team*\") {\n        $p.Range.Text = \"From April 2015 to July 2016, I lead my team in its projects and decision making. We remained autonomous within the organization, identifying problems and providing solutions wherever we saw value.\"\n        break\n    }\n}\n\n$d.Saved = $false\n"}
